$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old second header row ("Hiver"/"Eté"/"Année" units row); the data
# rows below shift up by one (old row 3 "Les Moyats" becomes row 2, etc.)
$ws.Rows.Item(2).Delete()

# Rewrite the (now single) header row with the new column headers
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("E1").ClearFormats()
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give F1:K1 a dedicated style (Arial 9, General format) distinct from the
# plain default used by A1:E1. We create a temporary named style to get Excel
# to register a new cell format (font applied, no explicit number format
# flag), apply it, then remove the named style again so only the underlying
# cell format remains - matching how the workbook's styles.xml ends up with
# one extra cellXfs entry but the same cellStyles/cellStyleXfs counts.
$tempStyle = $wb.Styles.Add("TempHeaderStyle")
$tempStyle.Font.Name = "Arial"
$tempStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TempHeaderStyle"
$wb.Styles.Item("TempHeaderStyle").Delete()

# Update the active selection to match the edited workbook
$null = $ws.Range("A2:K2").Select()
